# Update the single data row (row 2) of the API run-output sheet with the
# results of the latest simulation run (grid-connections population,
# parallelized main time-step loop).
#
# The source file stores every data cell as text (inline/shared string),
# even the numeric-looking ones, so we force each value to be written as
# text rather than letting Excel auto-detect it as a number. Formatting the
# cell as Text ("@") before assignment keeps the literal text (e.g. "1.0",
# "13.806") instead of Excel normalizing it to a number; ClearFormats()
# afterwards drops the temporary Text number-format again so the cell keeps
# the workbook's default (General) style, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("C2") "8761"        # timeStepsElapsed: 0 -> 8761
Set-TextValue $ws.Range("E2") "1.0"         # modelStartUpDuration_s: 0.439 -> 1.0
Set-TextValue $ws.Range("F2") "13.806"      # modelRunDuration_s: 0.007 -> 13.806
Set-TextValue $ws.Range("I2") "5739"        # nEnergyAssets: 5726 -> 5739
Set-TextValue $ws.Range("P2") "586.595"     # totalElectricityImported_MWh: 0.0 -> 586.595
Set-TextValue $ws.Range("Q2") "48258.307"   # totalElectricityExported_MWh: 0.0 -> 48258.307
